# Update computed "solidarity" statistics (mean / CI_low / CI_high) with
# final values from the re-run data preparation & render pipeline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0107402699915947
$ws.Range("C2").Value = -0.00289759134362644
$ws.Range("D2").Value = 0.0243781313268159

$ws.Range("B3").Value = 0.0165547590621099
$ws.Range("C3").Value = -0.00318621188945287
$ws.Range("D3").Value = 0.0362957300136726

$ws.Range("B4").Value = -0.00494056137424157
$ws.Range("C4").Value = -0.0512163291139684

$ws.Range("C5").Value = 0.00722864651991834
$ws.Range("D5").Value = 0.0890449926812002

$ws.Range("C6").Value = -0.022771040155234
$ws.Range("D6").Value = 0.0720264175854981

$ws.Range("D7").Value = 0.0470027632915175

$ws.Range("B8").Value = 0.0375599044291895
$ws.Range("C8").Value = -0.0220015874743622
$ws.Range("D8").Value = 0.0971213963327412

$ws.Range("B9").Value = 0.00518592478852574
$ws.Range("C9").Value = -0.0428998635668623

$ws.Range("C10").Value = -0.0757042568403408
$ws.Range("D10").Value = 0.0486503534558566

$ws.Range("B11").Value = 0.00955264618361334
$ws.Range("C11").Value = -0.0173418224556559
$ws.Range("D11").Value = 0.0364471148228826

$ws.Range("B12").Value = -0.0182522013113218
$ws.Range("C12").Value = -0.0630639428220745
$ws.Range("D12").Value = 0.0265595401994309

$ws.Range("B13").Value = -0.00920322016048767
$ws.Range("C13").Value = -0.0645346698358917

$ws.Range("B14").Value = 0.01779300719663
$ws.Range("C14").Value = -0.00836073172973517
$ws.Range("D14").Value = 0.0439467461229952
